$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1221.2
$ws.Range("I19").Value = 596
$ws.Range("K19").Value = 596
$ws.Range("M19").Value = -421
$ws.Range("H74").Value = 46139.07
$ws.Range("I74").Value = 61355.277
$ws.Range("K74").Value = 61355.277
$ws.Range("M74").Value = -60419.277
$ws.Range("H77").Value = 46139.07
$ws.Range("I77").Value = 61355.277
$ws.Range("K77").Value = 306776.385
$ws.Range("M77").Value = -302096.385
$ws.Range("H107").Value = 3352.5
$ws.Range("I107").Value = 2803.3333
$ws.Range("K107").Value = 2803.3333
$ws.Range("M107").Value = -883.3332999999998
$ws.Range("H137").Value = 2257
$ws.Range("I137").Value = 2535.3635
$ws.Range("J137").Value = 1916.7778
$ws.Range("K137").Value = 7606.0905
$ws.Range("L137").Value = 5750.3334
$ws.Range("M137").Value = -5056.0905
$ws.Range("N137").Value = -10850.3334
$ws.Range("H141").Value = 1155
$ws.Range("I141").Value = 1170
$ws.Range("K141").Value = 3510
$ws.Range("M141").Value = 1670

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 46669
$ws.Range("J24").Value = 46669
$ws.Range("L24").Value = 46669
$ws.Range("N24").Value = -47417
$ws.Range("H32").Value = 4050.9512
$ws.Range("I32").Value = 3175.4243
$ws.Range("K32").Value = 3175.4243
$ws.Range("M32").Value = -2888.4243
$ws.Range("H55").Value = 10126.333
$ws.Range("H61").Value = 4429.227
$ws.Range("I61").Value = 4270.684
$ws.Range("J61").Value = 5433.3335
$ws.Range("K61").Value = 4270.684
$ws.Range("L61").Value = 5433.3335
$ws.Range("M61").Value = -4058.684
$ws.Range("N61").Value = -5857.3335
$ws.Range("H74").Value = 1616
$ws.Range("I74").Value = 1277.8889
$ws.Range("K74").Value = 1277.8889
$ws.Range("M74").Value = -403.8888999999999
$ws.Range("H77").Value = 1616
$ws.Range("I77").Value = 1277.8889
$ws.Range("K77").Value = 6389.4445
$ws.Range("M77").Value = -2021.4445
$ws.Range("H100").Value = 46669
$ws.Range("J100").Value = 46669
$ws.Range("L100").Value = 46669
$ws.Range("N100").Value = -48833
$ws.Range("H102").Value = 2615.0476
$ws.Range("I102").Value = 2300.889
$ws.Range("K102").Value = 2300.889
$ws.Range("M102").Value = -678.8890000000001
$ws.Range("H132").Value = 2804.4092
$ws.Range("I132").Value = 1899.8334
$ws.Range("J132").Value = 6875
$ws.Range("K132").Value = 5699.5002
$ws.Range("L132").Value = 20625
$ws.Range("M132").Value = -3169.5002
$ws.Range("N132").Value = -25685
$ws.Range("H136").Value = 4429.227
$ws.Range("I136").Value = 4270.684
$ws.Range("J136").Value = 5433.3335
$ws.Range("K136").Value = 12812.052
$ws.Range("L136").Value = 16300.0005
$ws.Range("M136").Value = -10262.052
$ws.Range("N136").Value = -21400.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2891.6667
$ws.Range("I107").Value = 3126.4736
$ws.Range("K107").Value = 3126.4736
$ws.Range("M107").Value = -1206.4736
$ws.Range("H110").Value = 40702
$ws.Range("J110").Value = 40702
$ws.Range("L110").Value = 40702
$ws.Range("N110").Value = -48882

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 49563.227
$ws.Range("I31").Value = 60110.41
$ws.Range("J31").Value = 13702.8
$ws.Range("K31").Value = 60110.41
$ws.Range("L31").Value = 13702.8
$ws.Range("M31").Value = -59815.41
$ws.Range("N31").Value = -14292.8
$ws.Range("H34").Value = 49563.227
$ws.Range("I34").Value = 60110.41
$ws.Range("J34").Value = 13702.8
$ws.Range("K34").Value = 60110.41
$ws.Range("L34").Value = 13702.8
$ws.Range("M34").Value = -59908.41
$ws.Range("N34").Value = -14106.8
$ws.Range("H105").Value = 1498.091
$ws.Range("I105").Value = 1496.4445
$ws.Range("J105").Value = 1505.5
$ws.Range("K105").Value = 1496.4445
$ws.Range("L105").Value = 1505.5
$ws.Range("M105").Value = 250.5554999999999
$ws.Range("N105").Value = -4999.5
$ws.Range("H132").Value = 3264.4424
$ws.Range("I132").Value = 3213.1428
$ws.Range("K132").Value = 9639.428400000001
$ws.Range("M132").Value = -7109.428400000001
$ws.Range("H134").Value = 6935.1587
$ws.Range("I134").Value = 5009.1963
$ws.Range("K134").Value = 15027.5889
$ws.Range("M134").Value = -12492.5889

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 1160.4546
$ws.Range("I8").Value = 1160.4546
$ws.Range("K8").Value = 3481.3638
$ws.Range("M8").Value = -3342.3638
$ws.Range("H14").Value = 118371.06
$ws.Range("I14").Value = 118371.06
$ws.Range("K14").Value = 355113.18
$ws.Range("M14").Value = -354940.18
$ws.Range("H87").Value = 27162.375
$ws.Range("I87").Value = 16824.75
$ws.Range("K87").Value = 50474.25
$ws.Range("M87").Value = -49226.25
$ws.Range("H90").Value = 27162.375
$ws.Range("I90").Value = 16824.75
$ws.Range("K90").Value = 151422.75
$ws.Range("M90").Value = -145182.75
$ws.Range("H121").Value = 72048.82000000001
$ws.Range("I121").Value = 17599.5
$ws.Range("J121").Value = 79308.734
$ws.Range("K121").Value = 52798.5
$ws.Range("L121").Value = 237926.202
$ws.Range("M121").Value = -51488.5
$ws.Range("N121").Value = -240546.202

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2941.75
$ws.Range("I102").Value = 1184.5714
$ws.Range("K102").Value = 1184.5714
$ws.Range("M102").Value = 437.4286
$ws.Range("H132").Value = 560550.8
$ws.Range("J132").Value = 6401.5713
$ws.Range("L132").Value = 19204.7139
$ws.Range("N132").Value = -24264.7139

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4020.3635
$ws.Range("I132").Value = 3303
$ws.Range("J132").Value = 7248.5
$ws.Range("K132").Value = 9909
$ws.Range("L132").Value = 21745.5
$ws.Range("M132").Value = -7379
$ws.Range("N132").Value = -26805.5
$ws.Range("H136").Value = 3874.7407
$ws.Range("I136").Value = 3958.875
$ws.Range("J136").Value = 3201.6667
$ws.Range("K136").Value = 11876.625
$ws.Range("L136").Value = 9605.000100000001
$ws.Range("M136").Value = -9326.625
$ws.Range("N136").Value = -14705.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 37500
$ws.Range("J31").Value = 37500
$ws.Range("L31").Value = 37500
$ws.Range("N31").Value = -38196
$ws.Range("H122").Value = 2127.5151
$ws.Range("I122").Value = 2075.4644
$ws.Range("J122").Value = 2419
$ws.Range("K122").Value = 6226.3932
$ws.Range("L122").Value = 7257
$ws.Range("M122").Value = -3776.3932
$ws.Range("N122").Value = -12157

Write-Host "Applied all changes"